$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new log entry for "User story 2"
# Set in the same order the shared-string table picks up new unique
# strings (date, then subject, then description) so the produced
# sharedStrings.xml ordering matches Excel's own append order.
$ws.Range("C6").Value = "20-21.07.2024"
$ws.Range("A6").Value = "User story 2 completed"
$ws.Range("B6").Value = 12
$ws.Range("D6").Value = "Completed location, search, house feed, ordering houses"

# Update selection to match the author's final cursor position
$ws.Range("D7").Select() | Out-Null
